$wb = $excel.ActiveWorkbook

# Sheet 1: ip_address_list -> update E2 from 0 to 1
$wsList = $wb.Worksheets.Item("ip_address_list")
$wsList.Range("E2").Value = 1

# Sheet 2: ip_adress_fav_list -> add a new favorite row (row 1)
$wsFav = $wb.Worksheets.Item("ip_adress_fav_list")
$wsFav.Range("A1").NumberFormat = "@"
$wsFav.Range("A1").Value = "518"
$wsFav.Range("B1").Value = "192.168.1.241"
$wsFav.Range("C1").Value = "255.255.255.0"
$wsFav.Range("D1").Value = "pozngg"
$wsFav.Range("E1").Value = 1
